# Update the Singapore CONSTRUCTION_STANDARD workbook so that the
# STANDARD4 / STANDARD5 construction standards (rows 5 & 6) exist
# across every sheet, and the YEAR_START value used by the existing
# STANDARD1/STANDARD2/STANDARD3 rows on STANDARD_DEFINITION drops to 1000.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) STANDARD_DEFINITION: fix YEAR_START (col C) for rows 2-4,
#    then append two new rows (5 & 6) describing STANDARD4/STANDARD5.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("STANDARD_DEFINITION")

$ws1.Range("C2").Value = 1000
$ws1.Range("C3").Value = 1000
$ws1.Range("C4").Value = 1000

# Row 5 - new STANDARD4 entry (copy formatting + values from row 4, then tweak)
$ws1.Range("A4:D4").Copy()
$ws1.Range("A5:D5").PasteSpecial(-4122)
$ws1.Range("A4:D4").Copy()
$ws1.Range("A5:D5").PasteSpecial(-4163)
$ws1.Range("A5").Value = "STANDARD4"
$ws1.Range("B5").Value = "Concrete, Masonry and Rainscreens "
$ws1.Range("C5").Value = 1000
$ws1.Range("D5").Value = 2040

# Row 6 - duplicate row (source data keeps the STANDARD3 label here)
$ws1.Range("A4:D4").Copy()
$ws1.Range("A6:D6").PasteSpecial(-4122)
$ws1.Range("A4:D4").Copy()
$ws1.Range("A6:D6").PasteSpecial(-4163)
$ws1.Range("A6").Value = "STANDARD3"
$ws1.Range("B6").Value = "Concrete, Masonry and Rainscreens "
$ws1.Range("C6").Value = 1000
$ws1.Range("D6").Value = 2040

$ws1.Range("B10").Select()

# ---------------------------------------------------------------
# 2) ENVELOPE_ASSEMBLIES: append rows 5 & 6 (STANDARD4 / STANDARD5),
#    copying all assembly selections + values from row 4.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ENVELOPE_ASSEMBLIES")

$ws2.Range("A4:S4").Copy()
$ws2.Range("A5:S5").PasteSpecial(-4122)
$ws2.Range("A4:S4").Copy()
$ws2.Range("A5:S5").PasteSpecial(-4163)
$ws2.Range("A5").Value = "STANDARD4"

$ws2.Range("A4:S4").Copy()
$ws2.Range("A6:S6").PasteSpecial(-4122)
$ws2.Range("A4:S4").Copy()
$ws2.Range("A6:S6").PasteSpecial(-4163)
$ws2.Range("A6").Value = "STANDARD5"

$ws2.Range("B6:S6").Select()

# ---------------------------------------------------------------
# 3) HVAC_ASSEMBLIES: append rows 5 & 6 (STANDARD4 / STANDARD5).
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("HVAC_ASSEMBLIES")

$ws3.Range("A4:J4").Copy()
$ws3.Range("A5:J5").PasteSpecial(-4122)
$ws3.Range("A4:J4").Copy()
$ws3.Range("A5:J5").PasteSpecial(-4163)
$ws3.Range("A5").Value = "STANDARD4"

$ws3.Range("A4:J4").Copy()
$ws3.Range("A6:J6").PasteSpecial(-4122)
$ws3.Range("A4:J4").Copy()
$ws3.Range("A6:J6").PasteSpecial(-4163)
$ws3.Range("A6").Value = "STANDARD5"

$ws3.Range("B6:J6").Select()

# ---------------------------------------------------------------
# 4) SUPPLY_ASSEMBLIES: append rows 5 & 6 (STANDARD4 / STANDARD5).
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("SUPPLY_ASSEMBLIES")

$ws4.Range("A4:E4").Copy()
$ws4.Range("A5:E5").PasteSpecial(-4122)
$ws4.Range("A4:E4").Copy()
$ws4.Range("A5:E5").PasteSpecial(-4163)
$ws4.Range("A5").Value = "STANDARD4"

$ws4.Range("A4:E4").Copy()
$ws4.Range("A6:E6").PasteSpecial(-4122)
$ws4.Range("A4:E4").Copy()
$ws4.Range("A6:E6").PasteSpecial(-4163)
$ws4.Range("A6").Value = "STANDARD5"

$ws4.Range("B11").Select()

$ws1.Select()
